# Monte Carlo module added to PyDSS
# -----------------------------------------------------------------
# This script reproduces, via the Excel COM object model, the edit that
# added "Modbus" and "DNP3" as selectable Encoding values on the
# "Settings" sheet of the Socket Controller workbook, together with the
# knock-on changes to the "Controllers" sheet's data validation (the
# dropdown driving list grows from Settings!$A$3:$A$4 to
# Settings!$A$3:$A$6) and the various cosmetic UI state (selection,
# merged header, column widths) that moved along with it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Settings" sheet: insert "Modbus" and "DNP3" ahead of "ASN1"
# ---------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Before:                      After:
#  row3  None     / Server      row3  None     / Server
#  row4  ASN1     / Client      row4  Modbus   / Client
#                                row5  DNP3
#                                row6  ASN1
#
# Insert two fresh rows below the current "ASN1" row so it is pushed
# down to row 6, then relabel row 4 and fill in the two new rows.
$settings.Rows.Item(5).Resize(2, 1).EntireRow.Insert() | Out-Null

$settings.Range("A4").Value = "Modbus"
$settings.Range("A5").Value = "DNP3"
$settings.Range("A6").Value = "ASN1"

# The row-insert carries column B's formatting/value down onto the two
# new rows; only row 4 (Modbus/Client) should keep a "B" entry, so wipe
# the spurious B5/B6 cells completely (value + formatting).
$settings.Range("B5:B6").Clear()

# The header banner originally merged across A1:B1; it now spans the
# new helper column C as well.
$settings.Range("A1:B1").UnMerge()
$settings.Range("A1:C1").Merge()

# Approximate the new column C width (best-effort; exact OOXML "width"
# units depend on font metrics baked into the authoring copy of Excel).
$settings.Columns.Item(3).ColumnWidth = 9.5

# Restore the natural selection left behind in the authored workbook.
$settings.Range("A6").Select()

# ---------------------------------------------------------------
# 2) "Controllers" sheet: split the Encoding dropdown validation
# ---------------------------------------------------------------
$controllers = $wb.Worksheets.Item("Controllers")

# Previously E3:E1048576 (written as "E4:E1048576 E3") all shared one
# list validation against Settings!$A$3:$A$4. Now the header data row
# (E3) gets its own validation against the grown list
# Settings!$A$3:$A$6, while the rest of the column keeps referencing
# just Settings!$A$3:$A$4.
$controllers.Range("E3").Validation.Delete()
$controllers.Range("E3").Validation.Add(3, 1, 1, "=Settings!`$A`$3:`$A`$6")

$controllers.Range("E4:E1048576").Validation.Delete()
$controllers.Range("E4:E1048576").Validation.Add(3, 1, 1, "=Settings!`$A`$3:`$A`$4")

# Restore the natural selection left behind in the authored workbook.
$controllers.Activate()
$controllers.Range("H11").Select()

Write-Host "Monte Carlo module socket-controller edit applied"
